$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 57, pushing existing rows 57.. down by one.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new record's data.
$ws.Cells.Item(57, 1).Value = 8
$ws.Cells.Item(57, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(57, 3).Value = "Coquimbo"
$ws.Cells.Item(57, 4).Value = 44497
$ws.Cells.Item(57, 5).Value = 4
$ws.Cells.Item(57, 6).Value = 100114013
$ws.Cells.Item(57, 7).Value = "Zanahoria"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 600
$ws.Cells.Item(57, 11).Value = 6500
$ws.Cells.Item(57, 12).Value = 7000
$ws.Cells.Item(57, 13).Value = 6750
$ws.Cells.Item(57, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(57, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(57, 16).Value = 338
$ws.Cells.Item(57, 17).Value = 20
$ws.Cells.Item(57, 18).Value = "Hortaliza"
